$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update column C (Förändrad) for rows 2..506 from 45203 -> 45204
$ws.Range("C2:C506").Value = 45204

# Row 506 previously had no explicit row height; after the edit it gets the
# same default height that every other data row already carries.
$ws.Rows(506).RowHeight = 15

# 2) Append new row 507 with the new clearance notification record
$r = 507
$ws.Cells.Item($r, 1).Value = "A 47371-2023"
$ws.Cells.Item($r, 2).Value = 45202
$ws.Cells.Item($r, 2).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 3).Value = 45204
$ws.Cells.Item($r, 3).NumberFormat = "YYYY-MM-DD"
$ws.Cells.Item($r, 4).Value = "VÄRMLANDS LÄN"
$ws.Cells.Item($r, 5).Value = "KRISTINEHAMN"
$ws.Cells.Item($r, 7).Value = 10.1
$ws.Cells.Item($r, 8).Value = 0
$ws.Cells.Item($r, 9).Value = 0
$ws.Cells.Item($r, 10).Value = 0
$ws.Cells.Item($r, 11).Value = 0
$ws.Cells.Item($r, 12).Value = 0
$ws.Cells.Item($r, 13).Value = 0
$ws.Cells.Item($r, 14).Value = 0
$ws.Cells.Item($r, 15).Value = 0
$ws.Cells.Item($r, 16).Value = 0
$ws.Cells.Item($r, 17).Value = 0
$ws.Cells.Item($r, 18).WrapText = $true
